$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Loadcases")
$ws2 = $wb.Worksheets.Item("Envelopes")

# --- Sheet1 (Loadcases): new rows 14 & 15 (Explosion, Impact) ---
$ws1.Range("A14").Value = "Explosion"
$ws1.Range("B14").Value = 2

$ws1.Range("A15").Value = "Impact"
$ws1.Range("B15").Value = 2

# --- Sheet2 (Envelopes): new rows 18 & 19 (Explosion/Impact Envelope) ---
$ws2.Range("A18").Value = "Explosion Envelope"
$ws2.Range("B18").Value = "Explosion"
$ws2.Range("C18").Value = "Yes"

$ws2.Range("A19").Value = "Impact Envelope"
$ws2.Range("B19").Value = "Impact"
$ws2.Range("C19").Value = "Yes"

# --- Sheet1 (Loadcases): new row 16 (Seismic) ---
$ws1.Range("A16").Value = "Seismic"
$ws1.Range("B16").Value = 4
$ws1.Range("D16").Value = "Seismic Analysis"

# --- Sheet2 (Envelopes): new row 20 (Seismic Envelope) ---
$ws2.Range("A20").Value = "Seismic Envelope"
$ws2.Range("B20").Value = "Seismic"
$ws2.Range("C20").Value = "Yes"

# --- Sheet2: rename existing Action labels (frees old shared strings, appends new ones) ---
$ws2.Range("D2").Value = "00 Permanent Actions"

for ($r = 3; $r -le 17; $r++) {
    $ws2.Range("D$r").Value = "01 Variable Actions"
}

# --- Sheet2: Folder column for new rows ---
$ws2.Range("D18").Value = "02 Accidental Actions"
$ws2.Range("D19").Value = "02 Accidental Actions"
$ws2.Range("D20").Value = "03 Seismic Actions"

# --- Update selections / active sheet to match final state ---
$ws2.Activate()
$ws2.Range("I22").Select()
$ws1.Activate()
$ws1.Range("I18").Select()
